$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly data block (before former row 259),
# which shifts all existing data rows (259-360) down by two rows (to 261-362).
$ws.Range("A259:A260").EntireRow.Insert()

# Populate the newly inserted row 259 with the new week's "Primera" quality data.
$ws.Cells.Item(259, 1).Value = 11
$ws.Cells.Item(259, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(259, 3).Value = "Bíobío"
$ws.Cells.Item(259, 4).Value = 44924
$ws.Cells.Item(259, 5).Value = 8
$ws.Cells.Item(259, 6).Value = 100112008
$ws.Cells.Item(259, 7).Value = "Coliflor"
$ws.Cells.Item(259, 8).Value = "Sin especificar"
$ws.Cells.Item(259, 9).Value = "Primera"
$ws.Cells.Item(259, 10).Value = 2000
$ws.Cells.Item(259, 11).Value = 800
$ws.Cells.Item(259, 12).Value = 900
$ws.Cells.Item(259, 13).Value = 850
$ws.Cells.Item(259, 14).Value = "`$/unidad"
$ws.Cells.Item(259, 15).Value = "Región Metropolitana"
$ws.Cells.Item(259, 16).Value = 850
$ws.Cells.Item(259, 17).Value = 1
$ws.Cells.Item(259, 18).Value = "Hortaliza"

# Populate the newly inserted row 260 with the new week's "Segunda" quality data.
$ws.Cells.Item(260, 1).Value = 11
$ws.Cells.Item(260, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(260, 3).Value = "Bíobío"
$ws.Cells.Item(260, 4).Value = 44924
$ws.Cells.Item(260, 5).Value = 8
$ws.Cells.Item(260, 6).Value = 100112008
$ws.Cells.Item(260, 7).Value = "Coliflor"
$ws.Cells.Item(260, 8).Value = "Sin especificar"
$ws.Cells.Item(260, 9).Value = "Segunda"
$ws.Cells.Item(260, 10).Value = 1000
$ws.Cells.Item(260, 11).Value = 700
$ws.Cells.Item(260, 12).Value = 700
$ws.Cells.Item(260, 13).Value = 700
$ws.Cells.Item(260, 14).Value = "`$/unidad"
$ws.Cells.Item(260, 15).Value = "Región Metropolitana"
$ws.Cells.Item(260, 16).Value = 700
$ws.Cells.Item(260, 17).Value = 1
$ws.Cells.Item(260, 18).Value = "Hortaliza"
